$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 49,2
$data[0,0] = 0.9387094898216251
$data[0,1] = 3.328560344864989
$data[1,0] = 5.2515588033567
$data[1,1] = 7.056207800609791
$data[2,0] = 8.08938014359477
$data[2,1] = 11.09702766635007
$data[3,0] = 11.01824053167085
$data[3,1] = 14.66429942387451
$data[4,0] = 11.84569753089263
$data[4,1] = 18.21271944511849
$data[5,0] = 12.4495209675778
$data[5,1] = 22.5434771930722
$data[6,0] = 13.44748075756891
$data[6,1] = 26.39727248913607
$data[7,0] = 13.86066845125708
$data[7,1] = 30.3717169388179
$data[8,0] = 15.48518792122808
$data[8,1] = 33.8921025314127
$data[9,0] = 17.00679388538487
$data[9,1] = 38.54659097739153
$data[10,0] = 20.95140541880616
$data[10,1] = 42.24149132310051
$data[11,0] = 22.31574568595443
$data[11,1] = 46.41392157118887
$data[12,0] = 23.54100067858006
$data[12,1] = 50.77952037337192
$data[13,0] = 24.32186709528463
$data[13,1] = 54.48581403701866
$data[14,0] = 25.09131952196411
$data[14,1] = 59.16731212120687
$data[15,0] = 27.66190263377198
$data[15,1] = 63.27851700235214
$data[16,0] = 31.12760745988034
$data[16,1] = 67.66162291526179
$data[17,0] = 33.58960117990278
$data[17,1] = 71.48481870948902
$data[18,0] = 34.73103193208441
$data[18,1] = 74.69606751696463
$data[19,0] = 35.43751394950645
$data[19,1] = 78.63480631512766
$data[20,0] = 37.64885554961281
$data[20,1] = 82.7470089868835
$data[21,0] = 40.30079442421929
$data[21,1] = 86.17030902553435
$data[22,0] = 43.65326595969503
$data[22,1] = 89.91141423225264
$data[23,0] = 45.17561633605147
$data[23,1] = 93.26748394120146
$data[24,0] = 46.15529294915
$data[24,1] = 97.44676669850919
$data[25,0] = 49.91124275966857
$data[25,1] = 101.007297108561
$data[26,0] = 50.8176547077912
$data[26,1] = 104.950283133165
$data[27,0] = 54.45741191668857
$data[27,1] = 108.941864293514
$data[28,0] = 57.74499925511346
$data[28,1] = 113.4201246055219
$data[29,0] = 60.34465169762083
$data[29,1] = 116.8800621538798
$data[30,0] = 62.42278563121971
$data[30,1] = 120.4466170477735
$data[31,0] = 65.89985580528027
$data[31,1] = 124.6255167423677
$data[32,0] = 68.52561317154451
$data[32,1] = 128.1621175945332
$data[33,0] = 71.91799646663537
$data[33,1] = 132.5037106101475
$data[34,0] = 73.54985489992116
$data[34,1] = 136.1489835323165
$data[35,0] = 75.38510684341099
$data[35,1] = 139.7582609773989
$data[36,0] = 76.47177816999306
$data[36,1] = 143.463844443511
$data[37,0] = 77.89163740354344
$data[37,1] = 147.6389608462566
$data[38,0] = 79.5575859197231
$data[38,1] = 151.3519356546663
$data[39,0] = 80.72807075401498
$data[39,1] = 155.0888919349869
$data[40,0] = 84.88520353972142
$data[40,1] = 158.7408706707931
$data[41,0] = 86.47166134549553
$data[41,1] = 162.7205820157197
$data[42,0] = 89.55887749909935
$data[42,1] = 166.624619044602
$data[43,0] = 90.70465295740958
$data[43,1] = 170.6787376340409
$data[44,0] = 91.8868454220048
$data[44,1] = 174.6691552592767
$data[45,0] = 92.87421046636133
$data[45,1] = 178.8781241072097
$data[46,0] = 95.17754628350531
$data[46,1] = 182.493260175401
$data[47,0] = 97.46564754463925
$data[47,1] = 185.9804306053093
$data[48,0] = 99.91132497789691
$data[48,1] = 189.5984521261951

$ws.Range("B2:C50").Value = $data
